# Apply updated "powerx model" results to the sorted-results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - IPP369-Solar_1-Wind_1-ESS_2 (label unchanged, values updated)
$ws.Range("B2").Value = 4.135063959726197
$ws.Range("C2").Value = 29.9736666399491
$ws.Range("D2").Value = 18.75697678820473
$ws.Range("E2").Value = 9506.750456525891
$ws.Range("F2").Value = 10506.75045652589
$ws.Range("G2").Value = 596223854.8115311
$ws.Range("H2").Value = 65.00000000000007
$ws.Range("I2").Value = 62715.84150000006
$ws.Range("J2").Value = 30.00000000001278

# Row 3 - was IPP585-Solar_2-Wind_1-ESS_2, now IPP585-Solar_1-Wind_1-ESS_2
$ws.Range("A3").Value = "IPP585-Solar_1-Wind_1-ESS_2"
$ws.Range("B3").Value = 73.6457638928412
$ws.Range("C3").Value = -0
$ws.Range("D3").Value = 96.32362176882462
$ws.Range("E3").Value = 34335.64260191843
$ws.Range("F3").Value = 35335.64260191843
$ws.Range("G3").Value = 2153388719.222576
$ws.Range("H3").Value = 65.00000000000037
$ws.Range("I3").Value = 62715.84150000035
$ws.Range("J3").Value = 30.00000000003407

# Row 4 - was IPP585-Solar_2-Wind_2-ESS_2, now IPP585-Solar_1-Wind_2-ESS_2
$ws.Range("A4").Value = "IPP585-Solar_1-Wind_2-ESS_2"
$ws.Range("B4").Value = 73.6457638928412
$ws.Range("C4").Value = -0
$ws.Range("D4").Value = 96.32362176882462
$ws.Range("E4").Value = 34335.64260191843
$ws.Range("F4").Value = 35335.64260191843
$ws.Range("G4").Value = 2153388719.222576
$ws.Range("H4").Value = 65.00000000000037
$ws.Range("I4").Value = 62715.84150000035
$ws.Range("J4").Value = 30.00000000003407

# Row 5 - was IPP585-Solar_1-Wind_1-ESS_2, now IPP585-Solar_2-Wind_1-ESS_2
$ws.Range("A5").Value = "IPP585-Solar_2-Wind_1-ESS_2"
$ws.Range("B5").Value = 72.01408538931113
$ws.Range("C5").Value = -0
$ws.Range("D5").Value = 97.97904427143409
$ws.Range("E5").Value = 35101.2384516903
$ws.Range("F5").Value = 36101.2384516903
$ws.Range("G5").Value = 2201403707.189928
$ws.Range("H5").Value = 65.0000000000004
$ws.Range("I5").Value = 62715.84150000038
$ws.Range("J5").Value = 30.00000000001602

# Row 6 - was IPP585-Solar_1-Wind_2-ESS_2, now IPP585-Solar_2-Wind_2-ESS_2
$ws.Range("A6").Value = "IPP585-Solar_2-Wind_2-ESS_2"
$ws.Range("B6").Value = 72.01408538931113
$ws.Range("C6").Value = -0
$ws.Range("D6").Value = 97.97904427143409
$ws.Range("E6").Value = 35101.2384516903
$ws.Range("F6").Value = 36101.2384516903
$ws.Range("G6").Value = 2201403707.189928
$ws.Range("H6").Value = 65.0000000000004
$ws.Range("I6").Value = 62715.84150000038
$ws.Range("J6").Value = 30.00000000001602
